$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2965.2727
$ws.Range("I40").Value = 2105.5
$ws.Range("J40").Value = 3456.5715
$ws.Range("K40").Value = 2105.5
$ws.Range("L40").Value = 3456.5715
$ws.Range("M40").Value = -1930.5
$ws.Range("N40").Value = -3806.5715

$ws.Range("H61").Value = 132.77777
$ws.Range("I61").Value = 143.57143
$ws.Range("K61").Value = 430.71429
$ws.Range("M61").Value = -258.71429

$ws.Range("H134").Value = 166733140
$ws.Range("J134").Value = 166733140
$ws.Range("L134").Value = 166733140
$ws.Range("N134").Value = -166743280

$ws.Range("H138").Value = 2076.739
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2076.739
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 6230.217000000001
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -16510.217

$ws.Range("H141").Value = 9155.044
$ws.Range("I141").Value = 3958.6924
$ws.Range("J141").Value = 15910.3
$ws.Range("K141").Value = 11876.0772
$ws.Range("L141").Value = 47730.89999999999
$ws.Range("M141").Value = -6696.0772
$ws.Range("N141").Value = -58090.89999999999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8224.694
$ws.Range("J32").Value = 18720.4
$ws.Range("L32").Value = 18720.4
$ws.Range("N32").Value = -19294.4

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H61").Value = 2092.9167
$ws.Range("I61").Value = 2246.611
$ws.Range("J61").Value = 1939.2222
$ws.Range("K61").Value = 2246.611
$ws.Range("L61").Value = 1939.2222
$ws.Range("M61").Value = -2034.611
$ws.Range("N61").Value = -2363.2222

$ws.Range("H74").Value = 1676.4324
$ws.Range("I74").Value = 1581.92
$ws.Range("J74").Value = 1873.3334
$ws.Range("K74").Value = 1581.92
$ws.Range("L74").Value = 1873.3334
$ws.Range("M74").Value = -707.9200000000001
$ws.Range("N74").Value = -3621.3334

$ws.Range("H77").Value = 1676.4324
$ws.Range("I77").Value = 1581.92
$ws.Range("J77").Value = 1873.3334
$ws.Range("K77").Value = 7909.6
$ws.Range("L77").Value = 9366.666999999999
$ws.Range("M77").Value = -3541.6
$ws.Range("N77").Value = -18102.667

$ws.Range("H132").Value = 1449799.2
$ws.Range("I132").Value = 5297917
$ws.Range("J132").Value = 6755.125
$ws.Range("K132").Value = 15893751
$ws.Range("L132").Value = 20265.375
$ws.Range("M132").Value = -15891221
$ws.Range("N132").Value = -25325.375

$ws.Range("H136").Value = 2092.9167
$ws.Range("I136").Value = 2246.611
$ws.Range("J136").Value = 1939.2222
$ws.Range("K136").Value = 6739.833
$ws.Range("L136").Value = 5817.6666
$ws.Range("M136").Value = -4189.833
$ws.Range("N136").Value = -10917.6666


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 280.5
$ws.Range("I80").Value = 243.46666
$ws.Range("K80").Value = 243.46666
$ws.Range("M80").Value = 754.53334

$ws.Range("H83").Value = 280.5
$ws.Range("I83").Value = 243.46666
$ws.Range("K83").Value = 1217.3333
$ws.Range("M83").Value = 3774.6667

$ws.Range("H134").Value = 5872.9756
$ws.Range("I134").Value = 2803.1177
$ws.Range("J134").Value = 8047.4585
$ws.Range("K134").Value = 8409.3531
$ws.Range("L134").Value = 24142.3755
$ws.Range("M134").Value = -5874.3531
$ws.Range("N134").Value = -29212.3755


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4442
$ws.Range("J4").Value = 4786.2
$ws.Range("L4").Value = 4786.2
$ws.Range("N4").Value = -5010.2

$ws.Range("H31").Value = 3723.75
$ws.Range("I31").Value = 3482
$ws.Range("J31").Value = 3827.3572
$ws.Range("K31").Value = 3482
$ws.Range("L31").Value = 3827.3572
$ws.Range("M31").Value = -3187
$ws.Range("N31").Value = -4417.3572

$ws.Range("H34").Value = 3723.75
$ws.Range("I34").Value = 3482
$ws.Range("J34").Value = 3827.3572
$ws.Range("K34").Value = 3482
$ws.Range("L34").Value = 3827.3572
$ws.Range("M34").Value = -3280
$ws.Range("N34").Value = -4231.3572

$ws.Range("H58").Value = 2284.5264
$ws.Range("I58").Value = 1851.7778
$ws.Range("J58").Value = 3346.7273
$ws.Range("K58").Value = 1851.7778
$ws.Range("L58").Value = 3346.7273
$ws.Range("M58").Value = -1648.7778
$ws.Range("N58").Value = -3752.7273

$ws.Range("H132").Value = 2484.5454
$ws.Range("I132").Value = 1491.3846
$ws.Range("J132").Value = 3130.1
$ws.Range("K132").Value = 4474.1538
$ws.Range("L132").Value = 9390.299999999999
$ws.Range("M132").Value = -1944.1538
$ws.Range("N132").Value = -14450.3

$ws.Range("H136").Value = 2284.5264
$ws.Range("I136").Value = 1851.7778
$ws.Range("J136").Value = 3346.7273
$ws.Range("K136").Value = 5555.3334
$ws.Range("L136").Value = 10040.1819
$ws.Range("M136").Value = -3005.3334
$ws.Range("N136").Value = -15140.1819


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 601.7368
$ws.Range("I5").Value = 274.84616
$ws.Range("J5").Value = 1310
$ws.Range("K5").Value = 824.5384799999999
$ws.Range("L5").Value = 3930
$ws.Range("M5").Value = -712.5384799999999
$ws.Range("N5").Value = -4154

$ws.Range("H105").Value = 8965.916999999999
$ws.Range("J105").Value = 9324.091
$ws.Range("L105").Value = 27972.273
$ws.Range("N105").Value = -33214.273

$ws.Range("H122").Value = 1720.4242
$ws.Range("J122").Value = 1937.9259
$ws.Range("L122").Value = 17441.3331
$ws.Range("N122").Value = -22341.3331

$ws.Range("H135").Value = 601.7368
$ws.Range("I135").Value = 274.84616
$ws.Range("J135").Value = 1310
$ws.Range("K135").Value = 2473.61544
$ws.Range("L135").Value = 11790
$ws.Range("M135").Value = 61.38455999999996
$ws.Range("N135").Value = -16860


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1346646.4
$ws.Range("I132").Value = 3789918
$ws.Range("J132").Value = 2846.9
$ws.Range("K132").Value = 11369754
$ws.Range("L132").Value = 8540.700000000001
$ws.Range("M132").Value = -11367224
$ws.Range("N132").Value = -13600.7


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 698
$ws.Range("I22").Value = 372.5
$ws.Range("K22").Value = 372.5
$ws.Range("M22").Value = -77.5

$ws.Range("H27").Value = 698
$ws.Range("I27").Value = 372.5
$ws.Range("K27").Value = 372.5
$ws.Range("M27").Value = -265.5

$ws.Range("H122").Value = 9907.6
$ws.Range("I122").Value = 14838
$ws.Range("J122").Value = 4272.857
$ws.Range("K122").Value = 44514
$ws.Range("L122").Value = 12818.571
$ws.Range("M122").Value = -42064
$ws.Range("N122").Value = -17718.571

$ws.Range("H134").Value = 42414.5
$ws.Range("J134").Value = 42414.5
$ws.Range("L134").Value = 42414.5
$ws.Range("N134").Value = -52554.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1834.7709
$ws.Range("I132").Value = 1719.862
$ws.Range("J132").Value = 2010.1578
$ws.Range("K132").Value = 5159.586
$ws.Range("L132").Value = 6030.4734
$ws.Range("M132").Value = -2629.586
$ws.Range("N132").Value = -11090.4734

